$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B12").Value = "AAA"
$ws.Range("B15").Value = "BB"
$ws.Range("B18").Value = "BB"
$ws.Range("B30").Value = "AA"
$ws.Range("B39").Value = "BB"
$ws.Range("B51").Value = "B"
$ws.Range("B59").Value = "BB"
$ws.Range("B66").Value = "BBB"
$ws.Range("B71").Value = "BB"
$ws.Range("B76").Value = "BB"
$ws.Range("B87").Value = "A"
$ws.Range("B110").Value = "A"
$ws.Range("B112").Value = "A"
$ws.Range("B116").Value = "A"
$ws.Range("B118").Value = "A"
$ws.Range("B119").Value = "A"
$ws.Range("B122").Value = "BB"
$ws.Range("B129").Value = "A"
$ws.Range("B130").Value = "A"
$ws.Range("B134").Value = "A"
$ws.Range("B138").Value = "BB"
$ws.Range("B153").Value = "A"
$ws.Range("B155").Value = "BB"
$ws.Range("B158").Value = "BBB"
$ws.Range("B159").Value = "A"
$ws.Range("B161").Value = "AA"
$ws.Range("B163").Value = "AA"
$ws.Range("B180").Value = "AA"
$ws.Range("B181").Value = "AA"
$ws.Range("B182").Value = "BB"
$ws.Range("B184").Value = "B"
$ws.Range("B185").Value = "BB"
$ws.Range("B193").Value = "BB"
$ws.Range("B201").Value = "BB"
$ws.Range("B212").Value = "BB"
$ws.Range("B216").Value = "BB"
$ws.Range("B226").Value = "BBB"
$ws.Range("B232").Value = "AA"
$ws.Range("B233").Value = "BB"
$ws.Range("B237").Value = "B"
$ws.Range("B239").Value = "AAA"
$ws.Range("B265").Value = "BB"
$ws.Range("B270").Value = "BBB"
$ws.Range("B283").Value = "BB"
$ws.Range("B285").Value = "BB"
$ws.Range("B292").Value = "A"
$ws.Range("B297").Value = "BBB"
$ws.Range("B301").Value = "AA"
$ws.Range("B302").Value = "AAA"
$ws.Range("B306").Value = "A"
$ws.Range("B308").Value = "BBB"
$ws.Range("B311").Value = "A"
$ws.Range("B313").Value = "B"
$ws.Range("B323").Value = "B"
$ws.Range("B342").Value = "BB"
$ws.Range("B348").Value = "BB"
$ws.Range("B370").Value = "BB"
$ws.Range("B375").Value = "A"
$ws.Range("B387").Value = "BBB"
$ws.Range("B398").Value = "A"
